$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 292,
# shifting the former rows 292-320 down to 293-321.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new record's data.
$ws.Range("A292").Value = 7
$ws.Range("B292").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C292").Value = "Ñuble"
$ws.Range("D292").Value = 45106
$ws.Range("E292").Value = 16
$ws.Range("F292").Value = 100112032
$ws.Range("G292").Value = "Zapallo italiano"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 120
$ws.Range("K292").Value = 14000
$ws.Range("L292").Value = 15000
$ws.Range("M292").Value = 14500
$ws.Range("N292").Value = "$/caja 50 unidades"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 290
$ws.Range("Q292").Value = 50
$ws.Range("R292").Value = "Hortaliza"
